$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values (price and 1h volume % columns)
# Force text format on target cells so Excel does not auto-convert
# numeric-looking / percent-looking strings into Number/Percentage values,
# matching the original inlineStr (text) cell storage.
$cells = @{
    'D2' = '303.94'
    'E2' = '0.80%'
    'E3' = '11.42%'
    'D4' = '5.087'
    'E4' = '0.92%'
    'D5' = '0.07795'
    'E5' = '0.30%'
    'D6' = '2.275'
    'E6' = '-0.71%'
    'D7' = '8.086'
    'E7' = '1.71%'
    'D8' = '4.044'
    'E8' = '5.03%'
    'D9' = '0.9278'
    'E9' = '-0.45%'
    'D10' = '0.09149'
    'E10' = '-9.96%'
    'D11' = '0.1835'
    'E11' = '4.14%'
    'D12' = '0.08568'
    'E12' = '0.77%'
    'D13' = '0.03785'
    'E13' = '14.28%'
    'D14' = '0.09937'
    'E14' = '0.35%'
    'D15' = '0.001481'
    'E15' = '0.65%'
    'D16' = '0.005627'
    'E16' = '-1.31%'
    'D17' = '3.479'
    'E17' = '-0.42%'
    'D18' = '2.185'
    'E18' = '-0.30%'
    'E19' = '2.71%'
    'D20' = '0.1322'
    'E20' = '-1.50%'
    'D21' = '4.571'
    'E21' = '6.21%'
    'D22' = '0.2237'
    'E22' = '7.49%'
    'D23' = '0.04681'
    'E23' = '1.47%'
    'D24' = '0.001234'
    'E24' = '1.64%'
    'D25' = '0.004530'
    'E25' = '3.21%'
    'E26' = '1.17%'
    'E27' = '-19.71%'
    'D39' = '0.01780'
    'E39' = '3.66%'
    'D40' = '0.04728'
    'E40' = '-0.10%'
    'D41' = '0.007950'
    'E41' = '3.18%'
    'D42' = '0.1415'
    'E42' = '0.83%'
    'D43' = '0.007995'
    'E43' = '-17.95%'
    'D44' = '0.002304'
    'E44' = '11.55%'
    'D45' = '0.009621'
    'E45' = '-0.08%'
    'D46' = '0.00006221'
    'E46' = '2.73%'
    'D47' = '0.00000000752'
    'E47' = '1.14%'
    'D48' = '5.353'
    'E48' = '91.61%'
    'D49' = '0.002696'
    'E49' = '36.04%'
    'D50' = '0.00002104'
    'E50' = '1.14%'
    'D51' = '0.0002004'
    'E51' = '1.14%'
}

foreach ($addr in $cells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cells[$addr]
}
